# Update the EC (Estado de Cuenta) data table on Hoja1.
# The commit "Actualiza base de datos EC y agrega parte 1 de nuevos estado de
# cuenta" refreshes rows 16-25 (B:G) so the two workers' overdue periods
# (2308-2312) interleave correctly and the "2312" period amount reverts to
# 25333 while the rest use 46400.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "37617413"
$ws.Range("D16").Value = "YURISAN PATIÑO BOHORQUEZ"
$ws.Range("E16").Value = "2308"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1160000

# Row 17
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1050544690"
$ws.Range("D17").Value = "WILDER SANJUAN SERRANO"
$ws.Range("E17").Value = "2308"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1160000

# Row 18
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "37617413"
$ws.Range("D18").Value = "YURISAN PATIÑO BOHORQUEZ"
$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 46400
$ws.Range("G18").Value = 1160000

# Row 19
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1050544690"
$ws.Range("D19").Value = "WILDER SANJUAN SERRANO"
$ws.Range("E19").Value = "2309"
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1160000

# Row 20
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "37617413"
$ws.Range("D20").Value = "YURISAN PATIÑO BOHORQUEZ"
$ws.Range("E20").Value = "2310"
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1160000

# Row 21
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1050544690"
$ws.Range("D21").Value = "WILDER SANJUAN SERRANO"
$ws.Range("E21").Value = "2310"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1160000

# Row 22
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "37617413"
$ws.Range("D22").Value = "YURISAN PATIÑO BOHORQUEZ"
$ws.Range("E22").Value = "2311"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1160000

# Row 23
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1050544690"
$ws.Range("D23").Value = "WILDER SANJUAN SERRANO"
$ws.Range("E23").Value = "2311"
$ws.Range("F23").Value = 46400
$ws.Range("G23").Value = 1160000

# Row 24
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "37617413"
$ws.Range("D24").Value = "YURISAN PATIÑO BOHORQUEZ"
$ws.Range("E24").Value = "2312"
$ws.Range("F24").Value = 25333
$ws.Range("G24").Value = 1160000

# Row 25
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1050544690"
$ws.Range("D25").Value = "WILDER SANJUAN SERRANO"
$ws.Range("E25").Value = "2312"
$ws.Range("F25").Value = 25333
$ws.Range("G25").Value = 1160000
